# Generate Report for Handoff
# Inserts a new row (for file 70bbdfec-cd28-4f44-920e-fb9ddac3537f.md) above the
# existing a02783f8-... row on each of the three worksheets (Overview, zh-cn, de-de),
# pushing the original row down, extending each table, and fixing up the
# hyperlinks so they point at the correct files.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04b0623ccd1c951b7d867c8917073ac2d382c4af/e2e/"
$newFile = "70bbdfec-cd28-4f44-920e-fb9ddac3537f.md"
$oldFile = "a02783f8-0796-4b0c-bed2-1b4f1c5eb63a.md"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"  (columns A-G)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)

$ws.Rows.Item(2).Insert()
$lo.Resize($ws.Range("A1:G3"))

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = "e2e\" + $newFile
$ws.Range("C2").Value = ".md"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-31 14:50:42"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $baseUrl + $newFile, [System.Type]::Missing, [System.Type]::Missing, "e2e\" + $newFile)
$ws.Hyperlinks.Add($ws.Range("B3"), $baseUrl + $oldFile, [System.Type]::Missing, [System.Type]::Missing, "e2e\" + $oldFile)

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"  (columns A-P)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)

$ws.Rows.Item(2).Insert()
$lo.Resize($ws.Range("A1:P3"))

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "70bbdfec-cd28-4f44-920e-fb9ddac3537f.a434337ccbf9dece5f96e802f0a55ac5d3d6a12f.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-31 14:50:37"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + $newFile, [System.Type]::Missing, [System.Type]::Missing, $newFile)
$ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + $oldFile, [System.Type]::Missing, [System.Type]::Missing, $oldFile)

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"  (columns A-P)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)

$ws.Rows.Item(2).Insert()
$lo.Resize($ws.Range("A1:P3"))

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "70bbdfec-cd28-4f44-920e-fb9ddac3537f.a434337ccbf9dece5f96e802f0a55ac5d3d6a12f.de-de.xlf"
$ws.Range("H2").Value = "2016-08-31 14:50:42"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + $newFile, [System.Type]::Missing, [System.Type]::Missing, $newFile)
$ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + $oldFile, [System.Type]::Missing, [System.Type]::Missing, $oldFile)
